$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The date strings use a DD-MM-YYYY pattern; Excel's COM auto-detects some
# of these (where the day portion is <=12) as real dates and would silently
# convert them to date serials with a date number format applied. Force the
# cells to plain text first, then strip the resulting formatting afterwards
# so the saved cells keep their original (default) style while retaining
# the literal text value.
$dateRange = $ws.Range("A3:A21")
$dateRange.NumberFormat = "@"

$ws.Range("A3").Value = "28-07-2022"
$ws.Range("A4").Value = "01-08-2022"
$ws.Range("A5").Value = "04-08-2022"
$ws.Range("A6").Value = "08-08-2022"
$ws.Range("A7").Value = "11-08-2022"
$ws.Range("A8").Value = "15-08-2022"
$ws.Range("A9").Value = "18-08-2022"
$ws.Range("A10").Value = "22-08-2022"
$ws.Range("A11").Value = "25-08-2022"
$ws.Range("A12").Value = "29-08-2022"
$ws.Range("A13").Value = "01-09-2022"
$ws.Range("A14").Value = "05-09-2022"
$ws.Range("A15").Value = "08-09-2022"
$ws.Range("A16").Value = "12-09-2022"
$ws.Range("A17").Value = "15-09-2022"
$ws.Range("A18").Value = "19-09-2022"
$ws.Range("A19").Value = "22-09-2022"
$ws.Range("A20").Value = "26-09-2022"
$ws.Range("A21").Value = "29-09-2022"

# Remove the temporary text formatting so the cells keep no explicit style,
# matching the original workbook (only the inline text changes).
$dateRange.ClearFormats()

# Attendance count updates
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 0

$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("H6").Value = 0

$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 1
$ws.Range("H7").Value = 0

$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("H10").Value = 0

$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("H12").Value = 0
